# eqtl colocalisation for final credset July2025
#
# 1. Sheet2 ("var2gene" long table): a new eQTL/pops evidence row for
#    rs778801698 / RBM6 needs to be inserted right after the existing
#    "nearest_gene" row for the same variant (row 57), pushing the
#    following rows down by one.
# 2. Sheets 5, 6, 7, 8: drop the now-unused "sentinel_gtex" column (column O)
#    - it shifts every following column one to the left.
# 3. A brand-new Sheet9 is appended holding the PoPS gene-prioritisation
#    evidence for the same variant/gene.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet2 - insert the new "pops" evidence row for rs778801698 / RBM6
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 57 becomes a fresh copy of row 56 (same variant, same stats - only
# gene/evidence differ), then we overwrite those two columns.
$ws2.Rows.Item(57).Insert()
$ws2.Rows.Item(56).Copy()
$ws2.Rows.Item(57).PasteSpecial()
$excel.CutCopyMode = $false

$ws2.Range("E57").Value = "RBM6"
$ws2.Range("F57").Value = "pops"

# ---------------------------------------------------------------------
# 2) Sheets 5/6/7/8 - remove the "sentinel_gtex" column (column O)
# ---------------------------------------------------------------------
foreach ($sheetName in @("Sheet5", "Sheet6", "Sheet7", "Sheet8")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns.Item(15).Delete()
}

# ---------------------------------------------------------------------
# 3) Add Sheet9 (PoPS evidence) at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add($null, $lastSheet)
$ws9.Name = "Sheet9"

$headers = @(
    "locus","snpid","chr","posb37","posb38","a2","a1","PIP_average",
    "LOG_ODDS","se","eaf","pval","MAF","sentinel","sentinel_ubclung",
    "gene","ENSGID","PoPS_score","gene_rank","prioritized",
    "Feature1","Feature2","Feature3","Feature4","Feature5","Feature6",
    "Feature7","Feature8","Feature9","Feature10","gene_strand","evidence"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws9.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
}

$rowValues = @(
    "3_rs778801698_49524027_50524027",
    "rs778801698",
    3,
    50024027,
    49986594,
    "CA",
    "C",
    0.113624563121357,
    0.107858,
    0.0192236,
    0.565491,
    [double]"1.93263571048209e-08",
    0.434509,
    "3_50024027_CA_C",
    "3_50024027_C_CA",
    "RBM6",
    "ENSG00000004534",
    0.454640396848811,
    1,
    $true,
    "human_lung_diffexprs_genes_clusters_pre_def.36",
    "human_kidney_projected_pcaloadings.29",
    "human_colon_projected_icaloadings.6",
    "mouse_endothelium_projected_pcaloadings_clusters.276",
    "mouse_heart_control_projected_pcaloadings_clusters.98",
    "mouse_endothelium_projected_pcaloadings_clusters.350",
    "human_retina2_projected_pcaloadings_clusters.84",
    "human_muscle_projected_pcaloadings.95",
    "human_kidney3_projected_pcaloadings_clusters.31",
    "human_kidney3_average_expression.4",
    1,
    "pops"
)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws9.Cells.Item(2, $i + 1).Value = $rowValues[$i]
}

$ws9.Activate()
Write-Output "edit complete"
